$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bookings = @(
    @{ id = "SL-20251205-001"; created = "2025-12-05 12:01:55"; name = "Fahad Ahmed "; phone = "Fahad Ahmed "; tickets = 1; price = 175; total = 175; status = "pending"; intent = "70db35eb-9a82-4522-bdae-3f705a8ac0a6"; pstatus = "requires_payment_instrument" },
    @{ id = "SL-20251205-002"; created = "2025-12-05 12:02:04"; name = "Fahad Ahmed "; phone = "Fahad Ahmed "; tickets = 1; price = 175; total = 175; status = "pending"; intent = "3e0b0829-0a39-4f6a-b50d-799d1987355f"; pstatus = "requires_payment_instrument" },
    @{ id = "SL-20251205-003"; created = "2025-12-05 12:25:34"; name = "Fahad Ahmed "; phone = "Fahad Ahmed "; tickets = 1; price = 175; total = 175; status = "pending"; intent = "c65aa11b-2702-4d3b-b691-d0a2173d5350"; pstatus = "requires_payment_instrument" },
    @{ id = "SL-20251205-004"; created = "2025-12-05 12:34:42"; name = "Fahad Ahmed "; phone = "Fahad Ahmed "; tickets = 1; price = 175; total = 175; status = "pending"; intent = "ae5bcc5a-8027-4149-b15f-b746880a9c1a"; pstatus = "requires_payment_instrument" }
)

$row = 2
foreach ($b in $bookings) {
    $ws.Cells.Item($row, 1).Value = $b.id
    $ws.Cells.Item($row, 2).Value = $b.created
    $ws.Cells.Item($row, 3).Value = $b.name
    $ws.Cells.Item($row, 4).Value = $b.phone
    $ws.Cells.Item($row, 5).Value = $b.tickets
    $ws.Cells.Item($row, 6).Value = $b.price
    $ws.Cells.Item($row, 7).Value = $b.total
    $ws.Cells.Item($row, 8).Value = $b.status
    $ws.Cells.Item($row, 9).Value = $b.intent
    $ws.Cells.Item($row, 10).Value = $b.pstatus
    $ws.Cells.Item($row, 11).Value = "https://pay.ziina.com/payment_intent/$($b.intent)"
    $ws.Cells.Item($row, 12).Value = ""
    $row++
}
